# Repull data, push all data, mean calculation
# Updates the "dSF" column (column F) values on the sheet for the relevant rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    3  = 4
    4  = 2
    5  = -6
    6  = -1
    7  = 3
    8  = 1
    9  = -2
    10 = -5
    11 = 3
    13 = -1
    14 = -1
    16 = 6
    17 = 2
    18 = -1
    19 = -3
    20 = 7
    21 = 4
    22 = 1
    23 = 4
    24 = 4
    26 = -1
    27 = -5
    28 = -1
    29 = -1
    30 = -2
    31 = 1
    32 = 1
    33 = -1
    34 = -3
    36 = 3
    37 = -1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
